$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = " Cell A5"
$ws.Range("B4").Value = " Cell B4"
$ws.Range("B5").Value = " Cell B5"

[void]$ws.Range("F8").Select()
